# Update the "want to go" counts (column F) on each sheet.
# Each change bumps the existing count by a small amount, matching the
# regenerated data snapshot described in the commit message.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 465
$ws1.Range("F5").Value  = 82
$ws1.Range("F7").Value  = 635
$ws1.Range("F10").Value = 376
$ws1.Range("F13").Value = 343
$ws1.Range("F15").Value = 12860
$ws1.Range("F16").Value = 12846
$ws1.Range("F22").Value = 581
$ws1.Range("F23").Value = 2020
$ws1.Range("F25").Value = 12
$ws1.Range("F27").Value = 77
$ws1.Range("F28").Value = 256

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 11

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 170

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 170
$ws4.Range("F6").Value  = 465
$ws4.Range("F8").Value  = 82
$ws4.Range("F11").Value = 635
$ws4.Range("F15").Value = 376
$ws4.Range("F18").Value = 343
$ws4.Range("F21").Value = 12860
$ws4.Range("F22").Value = 12846
$ws4.Range("F28").Value = 581
$ws4.Range("F30").Value = 11
$ws4.Range("F31").Value = 2020
$ws4.Range("F33").Value = 12
$ws4.Range("F37").Value = 77
$ws4.Range("F38").Value = 256
